{"js": "const replacements = [\n  [\"2024-06-02 Sunday\", \"2024-06-03 Monday\"],\n  [\"78\u00d782=6396\", \"86\u00d762=5332\"],\n  [\"78\u00d751=3978\", \"60\u00d755=3300\"],\n  [\"71\u00d722=1562\", \"40\u00d784=3360\"],\n  [\"29\u00d718=522\", \"19\u00d770=1330\"],\n  [\"79\u00d727=2133\", \"30\u00d758=1740\"],\n  [\"99\u00d799=9801\", \"47\u00d746=2162\"],\n  [\"63\u00d716=1008\", \"20\u00d772=1440\"],\n  [\"65\u00d728=1820\", \"88\u00d746=4048\"],\n  [\"18\u00d767=1206\", \"62\u00d717=1054\"],\n  [\"81\u00d717=1377\", \"21\u00d782=1722\"],\n  [\"50\u00d748=2400\", \"14\u00d719=266\"],\n  [\"34\u00d734=1156\", \"31\u00d771=2201\"],\n  [\"74\u00d715=1110\", \"25\u00d780=2000\"],\n  [\"33\u00d756=1848\", \"58\u00d729=1682\"],\n  [\"75\u00d797=7275\", \"57\u00d799=5643\"],\n  [\"78\u00d730=2340\", \"84\u00d721=1764\"],\n  [\"76\u00d799=7524\", \"31\u00d760=1860\"],\n  [\"36\u00d760=2160\", \"31\u00d779=2449\"],\n  [\"32\u00d753=1696\", \"20\u00d754=1080\"],\n  [\"63\u00d743=2709\", \"38\u00d799=3762\"],\n  [\"50\u00d778=3900\", \"33\u00d732=1056\"],\n  [\"90\u00d722=1980\", \"26\u00d756=1456\"],\n  [\"12\u00d790=1080\", \"67\u00d777=5159\"],\n  [\"95\u00d796=9120\", \"88\u00d714=1232\"],\n  [\"95\u00d736=3420\", \"25\u00d783=2075\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph)\n$dateRange = $d.Paragraphs(1).Range\n$dateRange.MoveEnd(1, -1) | Out-Null\n$dateRange.Text = \"2024-06-03 Monday\"\n\n# Update each populated table cell (multiplication problems)\n$t = $d.Tables(1)\n$values = @(\n  @(1, @(\"86\u00d762=5332\", \"60\u00d755=3300\", \"40\u00d784=3360\", \"19\u00d770=1330\", \"30\u00d758=1740\")),\n  @(5, @(\"47\u00d746=2162\", \"20\u00d772=1440\", \"88\u00d746=4048\", \"62\u00d717=1054\", \"21\u00d782=1722\")),\n  @(10, @(\"14\u00d719=266\", \"31\u00d771=2201\", \"25\u00d780=2000\", \"58\u00d729=1682\", \"57\u00d799=5643\")),\n  @(15, @(\"84\u00d721=1764\", \"31\u00d760=1860\", \"31\u00d779=2449\", \"20\u00d754=1080\", \"38\u00d799=3762\")),\n  @(20, @(\"33\u00d732=1056\", \"26\u00d756=1456\", \"67\u00d777=5159\", \"88\u00d714=1232\", \"25\u00d783=2075\")),\n)\n\nforeach ($rowEntry in $values) {\n  $rowIndex = $rowEntry[0]\n  $cellTexts = $rowEntry[1]\n  for ($colIndex = 1; $colIndex -le $cellTexts.Count; $colIndex++) {\n    $cell = $t.Cell($rowIndex, $colIndex)\n    $cellRange = $cell.Range\n    $cellRange.MoveEnd(1, -1) | Out-Null\n    $cellRange.Text = $cellTexts[$colIndex - 1]\n  }\n}"}
